$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 currently holds the number 11038 (style index 9, numFmtId 3 "#,##0").
# The target turns it into the literal text "11.038", entered the same way a
# user would type it in Excel with a leading apostrophe (quote-prefix) so the
# existing number format is kept but the stored value becomes text.
$ws.Range("G2").Value = "'11.038"

# The author's selection moved to G3 after editing G2.
$ws.Range("G3").Select()
